$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value2 = 680.5
$ws.Range("I53").Value2 = 1037.4
$ws.Range("K53").Value2 = 1037.4
$ws.Range("M53").Value2 = -400.4000000000001
$ws.Range("H129").Value2 = 1418.5358
$ws.Range("I129").Value2 = 970.1539
$ws.Range("J129").Value2 = 1807.1333
$ws.Range("K129").Value2 = 2910.4617
$ws.Range("L129").Value2 = 5421.3999
$ws.Range("M129").Value2 = 2089.5383
$ws.Range("N129").Value2 = -15421.3999
$ws.Range("H138").Value2 = 1432663.9
$ws.Range("J138").Value2 = 2045548
$ws.Range("L138").Value2 = 6136644
$ws.Range("N138").Value2 = -6146924
$ws.Range("H140").Value2 = 105000
$ws.Range("J140").Value2 = 105000
$ws.Range("L140").Value2 = 105000
$ws.Range("N140").Value2 = -115360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2863666.8
$ws.Range("I32").Value2 = 3035252.5
$ws.Range("K32").Value2 = 3035252.5
$ws.Range("M32").Value2 = -3034965.5
$ws.Range("H39").Value2 = 0
$ws.Range("I39").Value2 = 0
$ws.Range("K39").Value2 = 0
$ws.Range("M39").ClearContents()
$ws.Range("H45").Value2 = 4709.0586
$ws.Range("I45").Value2 = 3291.75
$ws.Range("K45").Value2 = 3291.75
$ws.Range("M45").Value2 = -2914.75
$ws.Range("H51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("N51").ClearContents()
$ws.Range("H74").Value2 = 31443.514
$ws.Range("I74").Value2 = 43588.707
$ws.Range("K74").Value2 = 43588.707
$ws.Range("M74").Value2 = -42714.707
$ws.Range("H77").Value2 = 31443.514
$ws.Range("I77").Value2 = 43588.707
$ws.Range("K77").Value2 = 217943.535
$ws.Range("M77").Value2 = -213575.535
$ws.Range("H88").Value2 = 1686.091
$ws.Range("I88").Value2 = 0
$ws.Range("J88").Value2 = 1686.091
$ws.Range("K88").Value2 = 0
$ws.Range("L88").Value2 = 1686.091
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value2 = -2498.091
$ws.Range("H91").Value2 = 1686.091
$ws.Range("I91").Value2 = 0
$ws.Range("J91").Value2 = 1686.091
$ws.Range("K91").Value2 = 0
$ws.Range("L91").Value2 = 1686.091
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value2 = -4494.091
$ws.Range("H122").Value2 = 3399.0908
$ws.Range("I122").Value2 = 2673.111
$ws.Range("J122").Value2 = 6666
$ws.Range("K122").Value2 = 8019.333
$ws.Range("L122").Value2 = 19998
$ws.Range("M122").Value2 = -5569.333
$ws.Range("N122").Value2 = -24898
$ws.Range("H132").Value2 = 3498.4773
$ws.Range("I132").Value2 = 1173.7587
$ws.Range("J132").Value2 = 7992.933
$ws.Range("K132").Value2 = 3521.2761
$ws.Range("L132").Value2 = 23978.799
$ws.Range("M132").Value2 = -991.2761
$ws.Range("N132").Value2 = -29038.799

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value2 = 1120
$ws.Range("J5").Value2 = 1533.3334
$ws.Range("L5").Value2 = 1533.3334
$ws.Range("N5").Value2 = -1759.3334
$ws.Range("H62").Value2 = 0
$ws.Range("J62").Value2 = 0
$ws.Range("L62").Value2 = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("N65").ClearContents()
$ws.Range("H134").Value2 = 7357728.5
$ws.Range("I134").Value2 = 10872916
$ws.Range("J134").Value2 = 7790.4546
$ws.Range("K134").Value2 = 32618748
$ws.Range("L134").Value2 = 23371.3638
$ws.Range("M134").Value2 = -32616213
$ws.Range("N134").Value2 = -28441.3638

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 182
$ws.Range("I7").Value2 = 125.53846
$ws.Range("J7").Value2 = 263.55554
$ws.Range("K7").Value2 = 125.53846
$ws.Range("L7").Value2 = 263.55554
$ws.Range("M7").Value2 = -12.53846
$ws.Range("N7").Value2 = -489.55554
$ws.Range("H31").Value2 = 4443.857
$ws.Range("I31").Value2 = 1730.4642
$ws.Range("J31").Value2 = 7157.25
$ws.Range("K31").Value2 = 1730.4642
$ws.Range("L31").Value2 = 7157.25
$ws.Range("M31").Value2 = -1435.4642
$ws.Range("N31").Value2 = -7747.25
$ws.Range("H34").Value2 = 4443.857
$ws.Range("I34").Value2 = 1730.4642
$ws.Range("J34").Value2 = 7157.25
$ws.Range("K34").Value2 = 1730.4642
$ws.Range("L34").Value2 = 7157.25
$ws.Range("M34").Value2 = -1528.4642
$ws.Range("N34").Value2 = -7561.25
$ws.Range("H141").Value2 = 446641.72
$ws.Range("J141").Value2 = 446641.72
$ws.Range("L141").Value2 = 446641.72
$ws.Range("N141").Value2 = -457001.72

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 1248.875
$ws.Range("J5").Value2 = 2433.5
$ws.Range("L5").Value2 = 7300.5
$ws.Range("N5").Value2 = -7524.5
$ws.Range("H134").Value2 = 49258.137
$ws.Range("I134").Value2 = 53183.95
$ws.Range("K134").Value2 = 159551.85
$ws.Range("M134").Value2 = -154481.85
$ws.Range("H135").Value2 = 1248.875
$ws.Range("J135").Value2 = 2433.5
$ws.Range("L135").Value2 = 21901.5
$ws.Range("N135").Value2 = -26971.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value2 = 56538.5
$ws.Range("J62").Value2 = 43000
$ws.Range("L62").Value2 = 43000
$ws.Range("N62").Value2 = -44372
$ws.Range("H65").Value2 = 56538.5
$ws.Range("J65").Value2 = 43000
$ws.Range("L65").Value2 = 129000
$ws.Range("N65").Value2 = -135864
$ws.Range("H80").Value2 = 113677.555
$ws.Range("J80").Value2 = 168849.67
$ws.Range("L80").Value2 = 168849.67
$ws.Range("N80").Value2 = -170845.67
$ws.Range("H83").Value2 = 113677.555
$ws.Range("J83").Value2 = 168849.67
$ws.Range("L83").Value2 = 844248.3500000001
$ws.Range("N83").Value2 = -854232.3500000001
$ws.Range("H97").Value2 = 910.7778
$ws.Range("I97").Value2 = 1028.2727
$ws.Range("J97").Value2 = 726.1429000000001
$ws.Range("K97").Value2 = 1028.2727
$ws.Range("L97").Value2 = 726.1429000000001
$ws.Range("M97").Value2 = -532.2727
$ws.Range("N97").Value2 = -1718.1429
$ws.Range("H132").Value2 = 1677.5416
$ws.Range("I132").Value2 = 1000.8837
$ws.Range("J132").Value2 = 7496.8
$ws.Range("K132").Value2 = 3002.6511
$ws.Range("L132").Value2 = 22490.4
$ws.Range("M132").Value2 = -472.6511
$ws.Range("N132").Value2 = -27550.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 2271.9092
$ws.Range("I22").Value2 = 986.6667
$ws.Range("K22").Value2 = 986.6667
$ws.Range("M22").Value2 = -691.6667
$ws.Range("H27").Value2 = 2271.9092
$ws.Range("I27").Value2 = 986.6667
$ws.Range("K27").Value2 = 986.6667
$ws.Range("M27").Value2 = -879.6667
$ws.Range("H46").Value2 = 1985.625
$ws.Range("I46").Value2 = 1360.6666
$ws.Range("J46").Value2 = 3860.5
$ws.Range("K46").Value2 = 1360.6666
$ws.Range("L46").Value2 = 3860.5
$ws.Range("M46").Value2 = -1172.6666
$ws.Range("N46").Value2 = -4236.5
$ws.Range("H58").Value2 = 6250
$ws.Range("J58").Value2 = 8500
$ws.Range("L58").Value2 = 8500
$ws.Range("N58").Value2 = -9020
$ws.Range("H82").Value2 = 719405.5600000001
$ws.Range("J82").Value2 = 2122.6667
$ws.Range("L82").Value2 = 2122.6667
$ws.Range("N82").Value2 = -2844.6667
$ws.Range("H85").Value2 = 719405.5600000001
$ws.Range("J85").Value2 = 2122.6667
$ws.Range("L85").Value2 = 2122.6667
$ws.Range("N85").Value2 = -4618.6667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value2 = 59376
$ws.Range("J108").Value2 = 59376
$ws.Range("L108").Value2 = 59376
$ws.Range("N108").Value2 = -67056
$ws.Range("H113").Value2 = 1465.3334
$ws.Range("J113").Value2 = 1598.0769
$ws.Range("L113").Value2 = 4794.2307
$ws.Range("N113").Value2 = -9134.2307
$ws.Range("H132").Value2 = 3561.7407
$ws.Range("I132").Value2 = 3961.6365
$ws.Range("J132").Value2 = 2933.3333
$ws.Range("K132").Value2 = 11884.9095
$ws.Range("L132").Value2 = 8799.999899999999
$ws.Range("M132").Value2 = -9354.9095
$ws.Range("N132").Value2 = -13859.9999
